# Updates the crypto price (column D) and 1h volume/change (column E)
# figures produced by the scheduled GitHub Actions scrape.
#
# NOTE: several Price values (e.g. "0.4753", "315.36") are valid numeric
# literals, but the source data stores them as plain text (no thousands
# separator normalisation, fixed decimal widths, etc.). Assigning such a
# string directly to .Value would let Excel auto-convert it to a number.
# Prefixing with a literal apostrophe (a doubled '' inside a PowerShell
# single-quoted string) reproduces Excel's "quote prefix" input behavior,
# forcing the cell to stay text with the exact original string content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.465.13'
$ws.Range("E2").Value = '  +2.01%  '

$ws.Range("D3").Value = '1.843.23'
$ws.Range("E3").Value = '  +1.54%  '

$ws.Range("E4").Value = '  +1.31%  '

$ws.Range("D5").Value = '''315.36'
$ws.Range("E5").Value = '  +1.80%  '

$ws.Range("E6").Value = '  +1.12%  '

$ws.Range("D7").Value = '''0.4753'
$ws.Range("E7").Value = '  +1.42%  '

$ws.Range("D8").Value = '''0.3711'
$ws.Range("E8").Value = '  +0.57%  '

$ws.Range("D9").Value = '''0.07469'
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("D10").Value = '''0.8887'

$ws.Range("D11").Value = '''20.53'
$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("D12").Value = '1.845.57'
$ws.Range("E12").Value = '  +3.76%  '

$ws.Range("D13").Value = '''0.07390'
$ws.Range("E13").Value = '  +4.31%  '

$ws.Range("D14").Value = '''5.481'
$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").Value = '''93.59'
$ws.Range("E15").Value = '  +1.85%  '

$ws.Range("D16").Value = '''6.596'
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("E17").Value = '  +1.23%  '

$ws.Range("E18").Value = '  +1.62%  '

$ws.Range("E19").Value = '  +1.27%  '

$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("D21").Value = '27.477.44'
$ws.Range("E21").Value = '  +1.99%  '

$ws.Range("D22").Value = '''5.352'
$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("E23").Value = '  +1.21%  '

$ws.Range("D24").Value = '2.079.60'
$ws.Range("E24").Value = '  +2.77%  '

$ws.Range("D25").Value = '''1.902'
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("D26").Value = '''152.40'
$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("E27").Value = '  +1.63%  '

$ws.Range("D28").Value = '''2.173'
$ws.Range("E28").Value = '  -0.24%  '

$ws.Range("D29").Value = '''5.292'
$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("D30").Value = '''118.17'
$ws.Range("E30").Value = '  +1.81%  '

$ws.Range("D31").Value = '''0.08998'
$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("D32").Value = '''0.7607'
$ws.Range("E32").Value = '  -0.63%  '

$ws.Range("E33").Value = '  +1.33%  '

$ws.Range("D34").Value = '''4.572'
$ws.Range("E34").Value = '  +1.54%  '

$ws.Range("D35").Value = '''2.943'
$ws.Range("E35").Value = '  +0.75%  '

$ws.Range("E36").Value = '  +1.24%  '

$ws.Range("D37").Value = '''1.108'
$ws.Range("E37").Value = '  +1.99%  '

$ws.Range("E38").Value = '  +1.41%  '

$ws.Range("E39").Value = '  +0.43%  '

$ws.Range("E40").Value = '  +0.89%  '

$ws.Range("D41").Value = '''7.318'
$ws.Range("E41").Value = '  +0.87%  '

$ws.Range("D42").Value = '''0.5365'
$ws.Range("E42").Value = '  +0.38%  '

$ws.Range("D43").Value = '''2.373'
$ws.Range("E43").Value = '  +2.12%  '

$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("D45").Value = '''8.573'
$ws.Range("E45").Value = '  +1.47%  '

$ws.Range("D46").Value = '''0.4983'
$ws.Range("E46").Value = '  +1.21%  '

$ws.Range("D47").Value = '''10.68'

$ws.Range("E48").Value = '  +1.30%  '

$ws.Range("D49").Value = '''105.38'
$ws.Range("E49").Value = '  +2.23%  '

$ws.Range("E50").Value = '  +1.09%  '

$ws.Range("D51").Value = '''0.06323'
$ws.Range("E51").Value = '  +0.28%  '
